# The document has a "first page" header/footer pair and a "default"
# (all other pages) header/footer pair. Each one that carries a logo
# picture needs that picture's internal name swapped:
#   - first-page footer  (footer1.xml, docPr id="3"): image1.png -> image2.png
#   - default footer     (footer2.xml, docPr id="2"): image1.png -> image2.png
#   - first-page header  (header1.xml, docPr id="1"): image2.jpg -> image1.jpg

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# wdHeaderFooterFirstPage = 2, wdHeaderFooterPrimary = 1
$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2

# --- First page footer (footer1.xml) ---
$firstFooter = $section.Footers.Item($wdHeaderFooterFirstPage)
$firstFooterPic = $firstFooter.Range.InlineShapes.Item(1)
$firstFooterPic.Name = "image2.png"

# --- Default footer (footer2.xml) ---
$defaultFooter = $section.Footers.Item($wdHeaderFooterPrimary)
$defaultFooterPic = $defaultFooter.Range.InlineShapes.Item(1)
$defaultFooterPic.Name = "image2.png"

# --- First page header (header1.xml) ---
$firstHeader = $section.Headers.Item($wdHeaderFooterFirstPage)
$firstHeaderPic = $firstHeader.Range.InlineShapes.Item(1)
$firstHeaderPic.Name = "image1.jpg"

Write-Output "Renamed images in first-page footer, default footer and first-page header."
